$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), reusing H1's style (s="1")
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for columns I (I0) and J (IF), rows 2-24
$i0 = @(3, 6, 8, 8, 9, 8, 9, 9, 8, 6, 7, 8, 7, 5, 7, 8, 9, 8, 7, 5, 7, 8, 6)
$if = @(3, 7, 8, 8, 9, 8, 9, 9, 8, 7, 8, 9, 8, 7, 7, 9, 9, 8, 7, 5, 7, 8, 6)

for ($r = 0; $r -lt 23; $r++) {
    $row = $r + 2
    $ws.Cells.Item($row, 9).Value = $i0[$r]
    $ws.Cells.Item($row, 10).Value = $if[$r]
}
